# Add 2022-Q1 data
# 1. Create a new worksheet "2022-Q1" (positioned before "总计") by copying the
#    structurally-similar "2021-Q4" fund-holdings sheet (same headers/style),
#    then overwrite its data with the 2022-Q1 fund holdings.
# 2. Insert a new row into "总计" summarizing the 2022-Q1 quarter.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1: create the "2022-Q1" worksheet
# ---------------------------------------------------------------------------
$template = $wb.Worksheets.Item("2021-Q4")
$beforeSheet = $wb.Worksheets.Item("总计")
$template.Copy($beforeSheet)

# Sheet references can go stale (re-tracked by position) once the sheet
# collection changes, so re-fetch everything we need by name afterwards.
$newSheet = $wb.Worksheets.Item("2021-Q4 (2)")
$newSheet.Name = "2022-Q1"
$newSheet = $wb.Worksheets.Item("2022-Q1")

# The template ("2021-Q4") has 14 data rows (rows 2-15); we only need 9 data
# rows (rows 2-10), so drop the extra rows at the bottom.
$newSheet.Rows("11:15").Delete()

# Fund holdings data for 2022-Q1 (fund code, fund name, fund size, total
# stock position, position ratio, held market value (100M yuan), position rank)
$data = @(
  @("012260", "广发睿明优质企业混合型证券投资基金A", "12.09", "61.24", "2.90", "0.3506", 7),
  @("519673", "银河康乐股票",                         "2.31",  "92.35", "4.44", "0.1026", 6),
  @("011845", "博时周期优选混合型证券投资基金A",       "2.36",  "86.75", "2.83", "0.0668", 9),
  @("014212", "博时研究优享混合A",                     "2.30",  "60.01", "2.26", "0.0520", 10),
  @("012261", "广发睿明优质企业混合型证券投资基金C",   "0.84",  "61.24", "2.90", "0.0244", 7),
  @("217021", "招商优势企业混合",                       "0.36",  "69.72", "5.02", "0.0181", 6),
  @("014213", "博时研究优享混合C",                     "0.47",  "60.01", "2.26", "0.0106", 10),
  @("014157", "国泰君安创新医药混合",                   "0.31",  "79.44", "3.05", "0.0095", 3),
  @("011846", "博时周期优选混合型证券投资基金C",       "0.07",  "86.75", "2.83", "0.0020", 9)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rec = $data[$i]

    # Fund code -> text (preserve leading zeros)
    $c = $newSheet.Cells.Item($row, 2)
    $c.Value = "'" + $rec[0]
    $c.ClearFormats()

    # Fund name -> plain text
    $newSheet.Cells.Item($row, 3).Value = $rec[1]

    # Fund size / total stock position / position ratio / held market value
    # -> stored as text (mirrors the source data, keeps exact formatting such
    # as trailing zeros, e.g. "2.90")
    $c = $newSheet.Cells.Item($row, 4)
    $c.Value = "'" + $rec[2]
    $c.ClearFormats()

    $c = $newSheet.Cells.Item($row, 5)
    $c.Value = "'" + $rec[3]
    $c.ClearFormats()

    $c = $newSheet.Cells.Item($row, 6)
    $c.Value = "'" + $rec[4]
    $c.ClearFormats()

    $c = $newSheet.Cells.Item($row, 7)
    $c.Value = "'" + $rec[5]
    $c.ClearFormats()

    # Position rank -> number
    $newSheet.Cells.Item($row, 8).Value = $rec[6]
}

# ---------------------------------------------------------------------------
# Step 2: insert the summary row into "总计"
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Range("B2:D2").ClearFormats()

# Copy the bold/centered index-column style from row 3 (used by column A)
$totalSheet.Cells.Item(3, 1).Copy()
$totalSheet.Cells.Item(2, 1).PasteSpecial(-4122)

$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q1"
$totalSheet.Cells.Item(2, 3).Value = 9
$totalSheet.Cells.Item(2, 4).Value = 0.64

# Renumber the existing index column (A) for the rows that got shifted down
$totalSheet.Cells.Item(3, 1).Value = 1
$totalSheet.Cells.Item(4, 1).Value = 2
$totalSheet.Cells.Item(5, 1).Value = 3
$totalSheet.Cells.Item(6, 1).Value = 4

$excel.CutCopyMode = 0
